$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the discontinued "Volmax" row so the dataset shrinks from 80 to 79 rows
$ws.Rows.Item(74).Delete()

# Rewrite every data row (2-79) with the refreshed dataset / reordered items
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 'Augment'
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 'Augment 375 Tablet 18''s'
$ws.Cells.Item(2, 5).Value = '18''s'
$ws.Cells.Item(2, 6).Value = 38
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 'Augment'
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 'Augment 1.2g IV Injection 1''s'
$ws.Cells.Item(3, 5).Value = '1''s'
$ws.Cells.Item(3, 6).Value = 85
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = 'Augment'
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = 'Augment 625 Tablet 18''s'
$ws.Cells.Item(4, 5).Value = '18''s'
$ws.Cells.Item(4, 6).Value = 112
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = 'Augment'
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 4).Value = 'Augment 1gm Tablet 12''s'
$ws.Cells.Item(5, 5).Value = '12''s'
$ws.Cells.Item(5, 6).Value = 85
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = 'Augment'
$ws.Cells.Item(6, 3).Value = 5
$ws.Cells.Item(6, 4).Value = 'Augment 100ml PFS'
$ws.Cells.Item(6, 5).Value = '1''s'
$ws.Cells.Item(6, 6).Value = 65
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 2).Value = 'Biltin'
$ws.Cells.Item(7, 3).Value = 6
$ws.Cells.Item(7, 4).Value = 'Biltin 20mg Tablet 20''s'
$ws.Cells.Item(7, 5).Value = '20''s'
$ws.Cells.Item(7, 6).Value = 937
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = 'Bonflex'
$ws.Cells.Item(8, 3).Value = 7
$ws.Cells.Item(8, 4).Value = 'Bonflex FC Tablet'
$ws.Cells.Item(8, 5).Value = '40 ''s'
$ws.Cells.Item(8, 6).Value = 3
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = 'Desodin'
$ws.Cells.Item(9, 3).Value = 8
$ws.Cells.Item(9, 4).Value = 'Desodin 5mg Tablet'
$ws.Cells.Item(9, 5).Value = '50 ''s'
$ws.Cells.Item(9, 6).Value = 21
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = 'Dinafex'
$ws.Cells.Item(10, 3).Value = 12
$ws.Cells.Item(10, 4).Value = 'Dinafex 180mg FC Tablet 40''s'
$ws.Cells.Item(10, 5).Value = '40''s'
$ws.Cells.Item(10, 6).Value = 64
$ws.Cells.Item(11, 1).Value = 5
$ws.Cells.Item(11, 2).Value = 'Dinafex'
$ws.Cells.Item(11, 3).Value = 13
$ws.Cells.Item(11, 4).Value = 'Dinafex 120mg FC Tablet 40''s'
$ws.Cells.Item(11, 5).Value = '40''s'
$ws.Cells.Item(11, 6).Value = 291
$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = 'Dinafex'
$ws.Cells.Item(12, 3).Value = 14
$ws.Cells.Item(12, 4).Value = 'Dinafex 50ml Suspension'
$ws.Cells.Item(12, 5).Value = '50 ml'
$ws.Cells.Item(12, 6).Value = 774
$ws.Cells.Item(13, 1).Value = 5
$ws.Cells.Item(13, 2).Value = 'Dinafex'
$ws.Cells.Item(13, 3).Value = 16
$ws.Cells.Item(13, 4).Value = 'Dinafex 60mg FC Tablet 40''s'
$ws.Cells.Item(13, 5).Value = '40''s'
$ws.Cells.Item(13, 6).Value = 42
$ws.Cells.Item(14, 1).Value = 6
$ws.Cells.Item(14, 2).Value = 'Dorenta'
$ws.Cells.Item(14, 3).Value = 17
$ws.Cells.Item(14, 4).Value = 'Dorenta 100ml Syrup'
$ws.Cells.Item(14, 5).Value = '100 ml'
$ws.Cells.Item(14, 6).Value = 706
$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(15, 2).Value = 'Etorix'
$ws.Cells.Item(15, 3).Value = 20
$ws.Cells.Item(15, 4).Value = 'Etorix 60mg Tablet 50''s'
$ws.Cells.Item(15, 5).Value = '50''s'
$ws.Cells.Item(15, 6).Value = 575
$ws.Cells.Item(16, 1).Value = 7
$ws.Cells.Item(16, 2).Value = 'Etorix'
$ws.Cells.Item(16, 3).Value = 21
$ws.Cells.Item(16, 4).Value = 'Etorix 120mg Tablet 30''s'
$ws.Cells.Item(16, 5).Value = '30''s'
$ws.Cells.Item(16, 6).Value = 1310
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = 'Etorix'
$ws.Cells.Item(17, 3).Value = 22
$ws.Cells.Item(17, 4).Value = 'Etorix 90mg Tablet 40''s'
$ws.Cells.Item(17, 5).Value = '40''s'
$ws.Cells.Item(17, 6).Value = 1154
$ws.Cells.Item(18, 1).Value = 8
$ws.Cells.Item(18, 2).Value = 'Fenobac'
$ws.Cells.Item(18, 3).Value = 25
$ws.Cells.Item(18, 4).Value = 'Fenobac 5mg Tablet'
$ws.Cells.Item(18, 5).Value = '50 ''s'
$ws.Cells.Item(18, 6).Value = 20
$ws.Cells.Item(19, 1).Value = 8
$ws.Cells.Item(19, 2).Value = 'Fenobac'
$ws.Cells.Item(19, 3).Value = 27
$ws.Cells.Item(19, 4).Value = 'Fenobac 10mg Tablet'
$ws.Cells.Item(19, 5).Value = '30 ''s'
$ws.Cells.Item(19, 6).Value = 41
$ws.Cells.Item(20, 1).Value = 9
$ws.Cells.Item(20, 2).Value = 'Flucloxin'
$ws.Cells.Item(20, 3).Value = 30
$ws.Cells.Item(20, 4).Value = 'Flucloxin 100ml Dry Suspension'
$ws.Cells.Item(20, 5).Value = '100 ml'
$ws.Cells.Item(20, 6).Value = 2675
$ws.Cells.Item(21, 1).Value = 9
$ws.Cells.Item(21, 2).Value = 'Flucloxin'
$ws.Cells.Item(21, 3).Value = 31
$ws.Cells.Item(21, 4).Value = 'Flucloxin 500mg Capsule 40''s'
$ws.Cells.Item(21, 5).Value = '40''s'
$ws.Cells.Item(21, 6).Value = 1497
$ws.Cells.Item(22, 1).Value = 9
$ws.Cells.Item(22, 2).Value = 'Flucloxin'
$ws.Cells.Item(22, 3).Value = 32
$ws.Cells.Item(22, 4).Value = 'Flucloxin 250mg Capsule'
$ws.Cells.Item(22, 5).Value = '100 ''s'
$ws.Cells.Item(22, 6).Value = 156
$ws.Cells.Item(23, 1).Value = 9
$ws.Cells.Item(23, 2).Value = 'Flucloxin'
$ws.Cells.Item(23, 3).Value = 33
$ws.Cells.Item(23, 4).Value = 'Flucloxin 500mg IM/IV Injection'
$ws.Cells.Item(23, 5).Value = '1''s'
$ws.Cells.Item(23, 6).Value = 131
$ws.Cells.Item(24, 1).Value = 11
$ws.Cells.Item(24, 2).Value = 'Ketonic'
$ws.Cells.Item(24, 3).Value = 37
$ws.Cells.Item(24, 4).Value = 'Ketonic 60mg IM Injection'
$ws.Cells.Item(24, 5).Value = '1 ''s'
$ws.Cells.Item(24, 6).Value = 725
$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = 'Ketonic'
$ws.Cells.Item(25, 3).Value = 38
$ws.Cells.Item(25, 4).Value = 'Ketonic 10mg Tablet - 30''s'
$ws.Cells.Item(25, 5).Value = '30''s'
$ws.Cells.Item(25, 6).Value = 50
$ws.Cells.Item(26, 1).Value = 11
$ws.Cells.Item(26, 2).Value = 'Ketonic'
$ws.Cells.Item(26, 3).Value = 39
$ws.Cells.Item(26, 4).Value = 'Ketonic 30mg IM/IV Injection'
$ws.Cells.Item(26, 5).Value = '1 ''s'
$ws.Cells.Item(26, 6).Value = 747
$ws.Cells.Item(27, 1).Value = 13
$ws.Cells.Item(27, 2).Value = 'Levomax'
$ws.Cells.Item(27, 3).Value = 44
$ws.Cells.Item(27, 4).Value = 'Levomax 500mg Tablet - 20''s'
$ws.Cells.Item(27, 5).Value = '20 ''s'
$ws.Cells.Item(27, 6).Value = 59
$ws.Cells.Item(28, 1).Value = 13
$ws.Cells.Item(28, 2).Value = 'Levomax'
$ws.Cells.Item(28, 3).Value = 45
$ws.Cells.Item(28, 4).Value = 'Levomax 750mg Tablet - 10''s'
$ws.Cells.Item(28, 5).Value = '10 ''s'
$ws.Cells.Item(28, 6).Value = 2
$ws.Cells.Item(29, 1).Value = 14
$ws.Cells.Item(29, 2).Value = 'Lindamax'
$ws.Cells.Item(29, 3).Value = 46
$ws.Cells.Item(29, 4).Value = 'Lindamax 150mg Capsule'
$ws.Cells.Item(29, 5).Value = '30''s'
$ws.Cells.Item(29, 6).Value = 6
$ws.Cells.Item(30, 1).Value = 14
$ws.Cells.Item(30, 2).Value = 'Lindamax'
$ws.Cells.Item(30, 3).Value = 47
$ws.Cells.Item(30, 4).Value = 'Lindamax 300mg Capsule'
$ws.Cells.Item(30, 5).Value = '30''s'
$ws.Cells.Item(30, 6).Value = 25
$ws.Cells.Item(31, 1).Value = 14
$ws.Cells.Item(31, 2).Value = 'Lindamax'
$ws.Cells.Item(31, 3).Value = 48
$ws.Cells.Item(31, 4).Value = 'Lindamax Plus 10gm Gel'
$ws.Cells.Item(31, 5).Value = '10gm'
$ws.Cells.Item(31, 6).Value = 14
$ws.Cells.Item(32, 1).Value = 14
$ws.Cells.Item(32, 2).Value = 'Lindamax'
$ws.Cells.Item(32, 3).Value = 49
$ws.Cells.Item(32, 4).Value = 'Lindamax 25ml Lotion'
$ws.Cells.Item(32, 5).Value = '1''s'
$ws.Cells.Item(32, 6).Value = 7
$ws.Cells.Item(33, 1).Value = 15
$ws.Cells.Item(33, 2).Value = 'Mebidal'
$ws.Cells.Item(33, 3).Value = 50
$ws.Cells.Item(33, 4).Value = 'Mebidal Tablet'
$ws.Cells.Item(33, 5).Value = '200 ''s'
$ws.Cells.Item(33, 6).Value = 23
$ws.Cells.Item(34, 1).Value = 16
$ws.Cells.Item(34, 2).Value = 'Nabumet'
$ws.Cells.Item(34, 3).Value = 51
$ws.Cells.Item(34, 4).Value = 'Nabumet 750mg FC Tab 24''s'
$ws.Cells.Item(34, 5).Value = '24''s'
$ws.Cells.Item(34, 6).Value = 97
$ws.Cells.Item(35, 1).Value = 16
$ws.Cells.Item(35, 2).Value = 'Nabumet'
$ws.Cells.Item(35, 3).Value = 52
$ws.Cells.Item(35, 4).Value = 'Nabumet 500mg FC Tab 30''s'
$ws.Cells.Item(35, 5).Value = '30''s'
$ws.Cells.Item(35, 6).Value = 565
$ws.Cells.Item(36, 1).Value = 17
$ws.Cells.Item(36, 2).Value = 'Naprox'
$ws.Cells.Item(36, 3).Value = 53
$ws.Cells.Item(36, 4).Value = 'Naprox 15gm Gel'
$ws.Cells.Item(36, 5).Value = '15 gm'
$ws.Cells.Item(36, 6).Value = 12
$ws.Cells.Item(37, 1).Value = 17
$ws.Cells.Item(37, 2).Value = 'Naprox'
$ws.Cells.Item(37, 3).Value = 54
$ws.Cells.Item(37, 4).Value = 'Naprox 250mg Tablet'
$ws.Cells.Item(37, 5).Value = ' 50 ''s'
$ws.Cells.Item(37, 6).Value = 111
$ws.Cells.Item(38, 1).Value = 17
$ws.Cells.Item(38, 2).Value = 'Naprox'
$ws.Cells.Item(38, 3).Value = 55
$ws.Cells.Item(38, 4).Value = 'Naprox Plus 500mg Tablet - 36''s'
$ws.Cells.Item(38, 5).Value = '36''s'
$ws.Cells.Item(38, 6).Value = 357
$ws.Cells.Item(39, 1).Value = 17
$ws.Cells.Item(39, 2).Value = 'Naprox'
$ws.Cells.Item(39, 3).Value = 56
$ws.Cells.Item(39, 4).Value = 'Naprox 500mg Tablet'
$ws.Cells.Item(39, 5).Value = '50 ''s'
$ws.Cells.Item(39, 6).Value = 318
$ws.Cells.Item(40, 1).Value = 17
$ws.Cells.Item(40, 2).Value = 'Naprox'
$ws.Cells.Item(40, 3).Value = 57
$ws.Cells.Item(40, 4).Value = 'Naprox 50ml Suspension'
$ws.Cells.Item(40, 5).Value = '50 ml'
$ws.Cells.Item(40, 6).Value = 124
$ws.Cells.Item(41, 1).Value = 17
$ws.Cells.Item(41, 2).Value = 'Naprox'
$ws.Cells.Item(41, 3).Value = 59
$ws.Cells.Item(41, 4).Value = 'Naprox Plus 375mg Tablet - 30''s'
$ws.Cells.Item(41, 5).Value = '30 ''s'
$ws.Cells.Item(41, 6).Value = 106
$ws.Cells.Item(42, 1).Value = 18
$ws.Cells.Item(42, 2).Value = 'Ontin'
$ws.Cells.Item(42, 3).Value = 60
$ws.Cells.Item(42, 4).Value = 'Ontin 10mg Tablet'
$ws.Cells.Item(42, 5).Value = '100 ''s'
$ws.Cells.Item(42, 6).Value = 28
$ws.Cells.Item(43, 1).Value = 18
$ws.Cells.Item(43, 2).Value = 'Ontin'
$ws.Cells.Item(43, 3).Value = 61
$ws.Cells.Item(43, 4).Value = 'Ontin 60ml Syrup'
$ws.Cells.Item(43, 5).Value = '60 ml'
$ws.Cells.Item(43, 6).Value = 183
$ws.Cells.Item(44, 1).Value = 19
$ws.Cells.Item(44, 2).Value = 'Oradin'
$ws.Cells.Item(44, 3).Value = 62
$ws.Cells.Item(44, 4).Value = 'Oradin 60ml Suspension'
$ws.Cells.Item(44, 5).Value = '60 ml'
$ws.Cells.Item(44, 6).Value = 935
$ws.Cells.Item(45, 1).Value = 19
$ws.Cells.Item(45, 2).Value = 'Oradin'
$ws.Cells.Item(45, 3).Value = 63
$ws.Cells.Item(45, 4).Value = 'Oradin FT 10mg Tablet'
$ws.Cells.Item(45, 5).Value = '40 ''s'
$ws.Cells.Item(45, 6).Value = 16
$ws.Cells.Item(46, 1).Value = 19
$ws.Cells.Item(46, 2).Value = 'Oradin'
$ws.Cells.Item(46, 3).Value = 65
$ws.Cells.Item(46, 4).Value = 'Oradin 10mg Tablet'
$ws.Cells.Item(46, 5).Value = '100 ''s'
$ws.Cells.Item(46, 6).Value = 588
$ws.Cells.Item(47, 1).Value = 20
$ws.Cells.Item(47, 2).Value = 'Osticare'
$ws.Cells.Item(47, 3).Value = 66
$ws.Cells.Item(47, 4).Value = 'Osticare Tablet 30''s'
$ws.Cells.Item(47, 5).Value = '30''s'
$ws.Cells.Item(47, 6).Value = 156
$ws.Cells.Item(48, 1).Value = 20
$ws.Cells.Item(48, 2).Value = 'Osticare'
$ws.Cells.Item(48, 3).Value = 67
$ws.Cells.Item(48, 4).Value = 'Osticare FC Tab Container 30''s'
$ws.Cells.Item(48, 5).Value = '30''s'
$ws.Cells.Item(48, 6).Value = 17
$ws.Cells.Item(49, 1).Value = 21
$ws.Cells.Item(49, 2).Value = 'Paino'
$ws.Cells.Item(49, 3).Value = 69
$ws.Cells.Item(49, 4).Value = 'Paino 100mg Tablet'
$ws.Cells.Item(49, 5).Value = '100 ''s'
$ws.Cells.Item(49, 6).Value = 50
$ws.Cells.Item(50, 1).Value = 22
$ws.Cells.Item(50, 2).Value = 'Quinox'
$ws.Cells.Item(50, 3).Value = 70
$ws.Cells.Item(50, 4).Value = 'Quinox 500mg Tablet - 20''s'
$ws.Cells.Item(50, 5).Value = '20 ''s'
$ws.Cells.Item(50, 6).Value = 25
$ws.Cells.Item(51, 1).Value = 22
$ws.Cells.Item(51, 2).Value = 'Quinox'
$ws.Cells.Item(51, 3).Value = 71
$ws.Cells.Item(51, 4).Value = 'Quinox 500mg Tablet (40''s)'
$ws.Cells.Item(51, 5).Value = '40 ''s'
$ws.Cells.Item(51, 6).Value = 118
$ws.Cells.Item(52, 1).Value = 22
$ws.Cells.Item(52, 2).Value = 'Quinox'
$ws.Cells.Item(52, 3).Value = 72
$ws.Cells.Item(52, 4).Value = 'Quinox 750mg Tablet'
$ws.Cells.Item(52, 5).Value = '10 ''s'
$ws.Cells.Item(52, 6).Value = 9
$ws.Cells.Item(53, 1).Value = 22
$ws.Cells.Item(53, 2).Value = 'Quinox'
$ws.Cells.Item(53, 3).Value = 74
$ws.Cells.Item(53, 4).Value = 'Quinox 250mg Tablet'
$ws.Cells.Item(53, 5).Value = '30 ''s'
$ws.Cells.Item(53, 6).Value = 23
$ws.Cells.Item(54, 1).Value = 22
$ws.Cells.Item(54, 2).Value = 'Quinox'
$ws.Cells.Item(54, 3).Value = 75
$ws.Cells.Item(54, 4).Value = 'Quinox DS 60ml Pellets for Suspension'
$ws.Cells.Item(54, 5).Value = '60 ml'
$ws.Cells.Item(54, 6).Value = 600
$ws.Cells.Item(55, 1).Value = 23
$ws.Cells.Item(55, 2).Value = 'Rupaday'
$ws.Cells.Item(55, 3).Value = 76
$ws.Cells.Item(55, 4).Value = 'Rupaday 10mg Tablet 30''s'
$ws.Cells.Item(55, 5).Value = '30''s'
$ws.Cells.Item(55, 6).Value = 38
$ws.Cells.Item(56, 1).Value = 24
$ws.Cells.Item(56, 2).Value = 'Sk-Mox'
$ws.Cells.Item(56, 3).Value = 78
$ws.Cells.Item(56, 4).Value = 'Sk-Mox DS 100ml Dry Suspension'
$ws.Cells.Item(56, 5).Value = '100 ml'
$ws.Cells.Item(56, 6).Value = 40
$ws.Cells.Item(57, 1).Value = 24
$ws.Cells.Item(57, 2).Value = 'Sk-Mox'
$ws.Cells.Item(57, 3).Value = 79
$ws.Cells.Item(57, 4).Value = 'Sk-Mox 500mg Capsule 50''s'
$ws.Cells.Item(57, 5).Value = '50''s'
$ws.Cells.Item(57, 6).Value = 420
$ws.Cells.Item(58, 1).Value = 24
$ws.Cells.Item(58, 2).Value = 'Sk-Mox'
$ws.Cells.Item(58, 3).Value = 80
$ws.Cells.Item(58, 4).Value = 'Sk-Mox 15ml P/D'
$ws.Cells.Item(58, 5).Value = '15 ml'
$ws.Cells.Item(58, 6).Value = 75
$ws.Cells.Item(59, 1).Value = 24
$ws.Cells.Item(59, 2).Value = 'Sk-Mox'
$ws.Cells.Item(59, 3).Value = 81
$ws.Cells.Item(59, 4).Value = 'Sk-Mox 250mg Capsule'
$ws.Cells.Item(59, 5).Value = '100 ''s'
$ws.Cells.Item(59, 6).Value = 29
$ws.Cells.Item(60, 1).Value = 24
$ws.Cells.Item(60, 2).Value = 'Sk-Mox'
$ws.Cells.Item(60, 3).Value = 82
$ws.Cells.Item(60, 4).Value = 'Sk-Mox 100ml Dry Suspension'
$ws.Cells.Item(60, 5).Value = '100 ml'
$ws.Cells.Item(60, 6).Value = 469
$ws.Cells.Item(61, 1).Value = 24
$ws.Cells.Item(61, 2).Value = 'Sk-Mox'
$ws.Cells.Item(61, 3).Value = 83
$ws.Cells.Item(61, 4).Value = 'Sk-Mox 500mg Capsule'
$ws.Cells.Item(61, 5).Value = '48 ''s'
$ws.Cells.Item(61, 6).Value = 2
$ws.Cells.Item(62, 1).Value = 25
$ws.Cells.Item(62, 2).Value = 'Stiba'
$ws.Cells.Item(62, 3).Value = 84
$ws.Cells.Item(62, 4).Value = 'Stiba 10mg Tablet - 30''s'
$ws.Cells.Item(62, 5).Value = '30''s'
$ws.Cells.Item(62, 6).Value = 89
$ws.Cells.Item(63, 1).Value = 25
$ws.Cells.Item(63, 2).Value = 'Stiba'
$ws.Cells.Item(63, 3).Value = 85
$ws.Cells.Item(63, 4).Value = 'Stiba 30ml Syrup'
$ws.Cells.Item(63, 5).Value = '30ml'
$ws.Cells.Item(63, 6).Value = 10
$ws.Cells.Item(64, 1).Value = 26
$ws.Cells.Item(64, 2).Value = 'Sulidac'
$ws.Cells.Item(64, 3).Value = 86
$ws.Cells.Item(64, 4).Value = 'Sulidac 200 Tablet 20''s'
$ws.Cells.Item(64, 5).Value = '20''s'
$ws.Cells.Item(64, 6).Value = 101
$ws.Cells.Item(65, 1).Value = 26
$ws.Cells.Item(65, 2).Value = 'Sulidac'
$ws.Cells.Item(65, 3).Value = 87
$ws.Cells.Item(65, 4).Value = 'Sulidac 100mg Tablet 50''s'
$ws.Cells.Item(65, 5).Value = '50''s'
$ws.Cells.Item(65, 6).Value = 13
$ws.Cells.Item(66, 1).Value = 27
$ws.Cells.Item(66, 2).Value = 'Tenoxim'
$ws.Cells.Item(66, 3).Value = 88
$ws.Cells.Item(66, 4).Value = 'Tenoxim 20mg Tablet'
$ws.Cells.Item(66, 5).Value = '30''s'
$ws.Cells.Item(66, 6).Value = 66
$ws.Cells.Item(67, 1).Value = 28
$ws.Cells.Item(67, 2).Value = 'Timothy'
$ws.Cells.Item(67, 3).Value = 89
$ws.Cells.Item(67, 4).Value = 'Timothy 50mg Tablet'
$ws.Cells.Item(67, 5).Value = '50 ''s'
$ws.Cells.Item(67, 6).Value = 158
$ws.Cells.Item(68, 1).Value = 28
$ws.Cells.Item(68, 2).Value = 'Timothy'
$ws.Cells.Item(68, 3).Value = 90
$ws.Cells.Item(68, 4).Value = 'Timothy 5mg IM/IV'
$ws.Cells.Item(68, 5).Value = '5''s'
$ws.Cells.Item(68, 6).Value = 83
$ws.Cells.Item(69, 1).Value = 29
$ws.Cells.Item(69, 2).Value = 'Tojak'
$ws.Cells.Item(69, 3).Value = 91
$ws.Cells.Item(69, 4).Value = 'Tojak 5mg Tablet 10''s'
$ws.Cells.Item(69, 5).Value = '10''S'
$ws.Cells.Item(69, 6).Value = 212
$ws.Cells.Item(70, 1).Value = 30
$ws.Cells.Item(70, 2).Value = 'Toperin'
$ws.Cells.Item(70, 3).Value = 92
$ws.Cells.Item(70, 4).Value = 'Toperin 50mg Tablet - 60''s'
$ws.Cells.Item(70, 5).Value = '60 ''s'
$ws.Cells.Item(70, 6).Value = 31
$ws.Cells.Item(71, 1).Value = 31
$ws.Cells.Item(71, 2).Value = 'Toti'
$ws.Cells.Item(71, 3).Value = 93
$ws.Cells.Item(71, 4).Value = 'Toti 100ml Syrup'
$ws.Cells.Item(71, 5).Value = '100 ml'
$ws.Cells.Item(71, 6).Value = 3569
$ws.Cells.Item(72, 1).Value = 31
$ws.Cells.Item(72, 2).Value = 'Toti'
$ws.Cells.Item(72, 3).Value = 94
$ws.Cells.Item(72, 4).Value = 'Toti Tablet'
$ws.Cells.Item(72, 5).Value = '100 ''s'
$ws.Cells.Item(72, 6).Value = 74
$ws.Cells.Item(73, 1).Value = 32
$ws.Cells.Item(73, 2).Value = 'Visomox'
$ws.Cells.Item(73, 3).Value = 95
$ws.Cells.Item(73, 4).Value = 'Visomox 400mg FC Tablet 10''s'
$ws.Cells.Item(73, 5).Value = '10''s'
$ws.Cells.Item(73, 6).Value = 90
$ws.Cells.Item(74, 1).Value = 34
$ws.Cells.Item(74, 2).Value = 'Xenthol'
$ws.Cells.Item(74, 3).Value = 97
$ws.Cells.Item(74, 4).Value = 'Xenthol 30 Cream'
$ws.Cells.Item(74, 5).Value = '15gm'
$ws.Cells.Item(74, 6).Value = 9
$ws.Cells.Item(75, 1).Value = 35
$ws.Cells.Item(75, 2).Value = 'Zithrox'
$ws.Cells.Item(75, 3).Value = 98
$ws.Cells.Item(75, 4).Value = 'Zithrox 250mg Tablet - 12''s'
$ws.Cells.Item(75, 5).Value = '12''s'
$ws.Cells.Item(75, 6).Value = 48
$ws.Cells.Item(76, 1).Value = 35
$ws.Cells.Item(76, 2).Value = 'Zithrox'
$ws.Cells.Item(76, 3).Value = 99
$ws.Cells.Item(76, 4).Value = 'Zithrox 500mg Tablet - 12''s'
$ws.Cells.Item(76, 5).Value = '12 ''s'
$ws.Cells.Item(76, 6).Value = 679
$ws.Cells.Item(77, 1).Value = 35
$ws.Cells.Item(77, 2).Value = 'Zithrox'
$ws.Cells.Item(77, 3).Value = 101
$ws.Cells.Item(77, 4).Value = 'Zithrox 35ml Dry Suspension'
$ws.Cells.Item(77, 5).Value = '35ml'
$ws.Cells.Item(77, 6).Value = 825
$ws.Cells.Item(78, 1).Value = 35
$ws.Cells.Item(78, 2).Value = 'Zithrox'
$ws.Cells.Item(78, 3).Value = 103
$ws.Cells.Item(78, 4).Value = 'Zithrox 20ml Powder for Suspension'
$ws.Cells.Item(78, 5).Value = '20ml'
$ws.Cells.Item(78, 6).Value = 287
$ws.Cells.Item(79, 1).Value = 35
$ws.Cells.Item(79, 2).Value = 'Zithrox'
$ws.Cells.Item(79, 3).Value = 106
$ws.Cells.Item(79, 4).Value = 'Zithrox 50ml Powder for Suspension'
$ws.Cells.Item(79, 5).Value = '50ml'
$ws.Cells.Item(79, 6).Value = 88
